$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 5 duplicates row 4's values/formatting (new "checkSavedSearch" test case
# case shares the same flight/traveller data as the existing "checkFilters" row).
$ws.Range("A4").Copy($ws.Range("A5")) | Out-Null
$ws.Range("B4").Copy($ws.Range("B5")) | Out-Null
$ws.Range("C4").Copy($ws.Range("C5")) | Out-Null
$ws.Range("D4").Copy($ws.Range("D5")) | Out-Null
$ws.Range("E4").Copy($ws.Range("E5")) | Out-Null
$ws.Range("F4").Copy($ws.Range("F5")) | Out-Null
$ws.Range("G4").Copy($ws.Range("G5")) | Out-Null
$ws.Range("H4").Copy($ws.Range("H5")) | Out-Null
$ws.Range("I4").Copy($ws.Range("I5")) | Out-Null
$ws.Range("J4").Copy($ws.Range("J5")) | Out-Null
$ws.Range("K4").Copy($ws.Range("K5")) | Out-Null
$ws.Range("L4").Copy($ws.Range("L5")) | Out-Null

$ws.Range("A5").Value = "checkSavedSearch"

# Column A now needs to fit the widest TestCaseName value ("checkSavedSearch").
$ws.Columns.Item(1).ColumnWidth = 17.42578125

# Active selection moves to the last populated cell in the new row.
$ws.Range("L5").Select() | Out-Null
